$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextCell 2 4 '35.258.09'  # D2: was '35.317.85'
Set-TextCell 2 5 '  -0.19%  '  # E2: was '  -0.09%  '
Set-TextCell 3 4 '1.907.64'  # D3: was '1.912.67'
Set-TextCell 3 5 '  -0.04%  '  # E3: was '  +0.20%  '
Set-TextCell 4 5 '  +0.33%  '  # E4: was '  +0.04%  '
Set-TextCell 5 2 'BNB'  # B5: was 'XRP'
Set-TextCell 5 3 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'  # C5: was 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextCell 5 4 '255.08'  # D5: was '0.725'
Set-TextCell 5 5 '  +3.37%  '  # E5: was '  +8.62%  '
Set-TextCell 6 2 'XRP'  # B6: was 'BNB'
Set-TextCell 6 3 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'  # C6: was 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextCell 6 4 '0.694'  # D6: was '255.30'
Set-TextCell 6 5 '  +1.30%  '  # E6: was '  +3.61%  '
Set-TextCell 7 5 '  +0.35%  '  # E7: was '  +0.09%  '
Set-TextCell 8 4 '42.15'  # D8: was '42.38'
Set-TextCell 8 5 '  +1.13%  '  # E8: was '  +1.87%  '
Set-TextCell 9 4 '0.364'  # D9: was '0.366'
Set-TextCell 9 5 '  +4.10%  '  # E9: was '  +5.84%  '
Set-TextCell 10 4 '52.94'  # D10: was '53.27'
Set-TextCell 10 5 '  +1.13%  '  # E10: was '  +0.38%  '
Set-TextCell 11 4 '0.0755'  # D11: was '0.0771'
Set-TextCell 11 5 '  +3.74%  '  # E11: was '  +7.29%  '
Set-TextCell 12 4 '0.0988'  # D12: was '0.0987'
Set-TextCell 12 5 '  -0.56%  '  # E12: was '  -0.48%  '
Set-TextCell 13 2 'WrappedliquidstakedEther2.0'  # B13: was 'Chainlink'
Set-TextCell 13 3 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'  # C13: was 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell 13 4 '2.189.29'  # D13: was '13.10'
Set-TextCell 13 5 '  +0.12%  '  # E13: was '  +6.34%  '
Set-TextCell 14 2 'Chainlink'  # B14: was 'WrappedliquidstakedEther2.0'
Set-TextCell 14 3 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'  # C14: was 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextCell 14 4 '12.92'  # D14: was '2.189.55'
Set-TextCell 14 5 '  +3.87%  '  # E14: was '  +0.19%  '
Set-TextCell 15 4 '0.733'  # D15: was '0.739'
Set-TextCell 15 5 '  +4.26%  '  # E15: was '  +5.70%  '
Set-TextCell 16 4 '4.96'  # D16: was '5.01'
Set-TextCell 16 5 '  +2.60%  '  # E16: was '  +3.96%  '
Set-TextCell 17 4 '1.921.82'  # D17: was '1.904.49'
Set-TextCell 17 5 '  +0.74%  '  # E17: was '  -0.30%  '
Set-TextCell 18 4 '35.316.28'  # D18: was '35.317.75'
Set-TextCell 18 5 '  -0.04%  '  # E18: was '  -0.14%  '
Set-TextCell 19 4 '74.41'  # D19: was '75.00'
Set-TextCell 19 5 '  +2.56%  '  # E19: was '  +3.93%  '
Set-TextCell 20 4 '0.0₃0835'  # D20: was '0.0₃0849'
Set-TextCell 20 5 '  +1.44%  '  # E20: was '  +3.51%  '
Set-TextCell 21 4 '243.48'  # D21: was '245.61'
Set-TextCell 21 5 '  +1.25%  '  # E21: was '  +1.96%  '
Set-TextCell 22 4 '13.04'  # D22: was '13.17'
Set-TextCell 22 5 '  +4.29%  '  # E22: was '  +5.57%  '
Set-TextCell 23 4 '5.05'  # D23: was '5.16'
Set-TextCell 23 5 '  +3.89%  '  # E23: was '  +6.99%  '
Set-TextCell 24 5 '  +0.33%  '  # E24: was '  +0.04%  '
Set-TextCell 25 4 '2.45'  # D25: was '2.46'
Set-TextCell 25 5 '  +6.89%  '  # E25: was '  +7.32%  '
Set-TextCell 26 4 '2.35'  # D26: was '2.39'
Set-TextCell 26 5 '  -3.37%  '  # E26: was '  -0.37%  '
Set-TextCell 27 4 '166.64'  # D27: was '166.89'
Set-TextCell 27 5 '  -2.31%  '  # E27: was '  -2.40%  '
Set-TextCell 28 4 '8.63'  # D28: was '8.80'
Set-TextCell 28 5 '  +0.62%  '  # E28: was '  +4.11%  '
Set-TextCell 29 4 '18.58'  # D29: was '18.84'
Set-TextCell 29 5 '  +0.60%  '  # E29: was '  +2.69%  '
Set-TextCell 30 4 '0.130'  # D30: was '0.132'
Set-TextCell 30 5 '  +0.20%  '  # E30: was '  +4.23%  '
Set-TextCell 31 4 '4.127.68'  # D31: was '4.128.81'
Set-TextCell 31 5 '  -0.57%  '  # E31: was '  -0.54%  '
Set-TextCell 32 2 'WEMIXToken'  # B32: was 'TrustWalletToken'
Set-TextCell 32 3 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'  # C32: was 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell 32 4 '2.01'  # D32: was '1.69'
Set-TextCell 32 5 '  +13.58%  '  # E32: was '  +27.03%  '
Set-TextCell 33 4 '4.33'  # D33: was '4.36'
Set-TextCell 33 5 '  +3.12%  '  # E33: was '  +5.20%  '
Set-TextCell 34 2 'TrustWalletToken'  # B34: was 'WEMIXToken'
Set-TextCell 34 3 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'  # C34: was 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextCell 34 4 '1.65'  # D34: was '2.00'
Set-TextCell 34 5 '  +23.78%  '  # E34: was '  +15.32%  '
Set-TextCell 35 4 '0.0581'  # D35: was '0.0593'
Set-TextCell 35 5 '  +2.12%  '  # E35: was '  +4.92%  '
Set-TextCell 36 4 '4.21'  # D36: was '4.28'
Set-TextCell 36 5 '  +2.25%  '  # E36: was '  +4.55%  '
Set-TextCell 37 5 '  +0.35%  '  # E37: was '  +0.05%  '
Set-TextCell 38 4 '0.872'  # D38: was '0.929'
Set-TextCell 38 5 '  -11.52%  '  # E38: was '  -1.86%  '
Set-TextCell 39 4 '2.01'  # D39: was '2.04'
Set-TextCell 39 5 '  -1.82%  '  # E39: was '  +0.00%  '
Set-TextCell 40 4 '98.81'  # D40: was '100.08'
Set-TextCell 40 5 '  +9.66%  '  # E40: was '  +11.14%  '
Set-TextCell 41 4 '0.0216'  # D41: was '0.0220'
Set-TextCell 41 5 '  +3.73%  '  # E41: was '  +5.88%  '
Set-TextCell 42 2 'InjectiveProtocol'  # B42: was 'ARBITRUM'
Set-TextCell 42 3 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'  # C42: was 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell 42 4 '16.99'  # D42: was '1.13'
Set-TextCell 42 5 '  +3.71%  '  # E42: was '  +2.69%  '
Set-TextCell 43 2 'ARBITRUM'  # B43: was 'InjectiveProtocol'
Set-TextCell 43 3 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'  # C43: was 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell 43 4 '1.12'  # D43: was '17.03'
Set-TextCell 43 5 '  +0.47%  '  # E43: was '  +4.70%  '
Set-TextCell 44 4 '0.0644'  # D44: was '0.0650'
Set-TextCell 44 5 '  -3.44%  '  # E44: was '  -0.09%  '
Set-TextCell 45 2 'Maker'  # B45: was 'RenderToken'
Set-TextCell 45 3 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'  # C45: was 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 45 4 '1.334.38'  # D45: was '2.47'
Set-TextCell 45 5 '  -0.62%  '  # E45: was '  +3.61%  '
Set-TextCell 46 2 'RenderToken'  # B46: was 'Maker'
Set-TextCell 46 3 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'  # C46: was 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell 46 4 '2.44'  # D46: was '1.342.23'
Set-TextCell 46 5 '  +2.09%  '  # E46: was '  +0.24%  '
Set-TextCell 47 5 '  +1.66%  '  # E47: was '  +1.37%  '
Set-TextCell 48 2 'MXToken'  # B48: was 'FraxShare'
Set-TextCell 48 3 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'  # C48: was 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 48 4 '2.76'  # D48: was '6.73'
Set-TextCell 48 5 '  -1.46%  '  # E48: was '  +3.40%  '
Set-TextCell 49 2 'FraxShare'  # B49: was 'MXToken'
Set-TextCell 49 3 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'  # C49: was 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell 49 4 '6.61'  # D49: was '2.76'
Set-TextCell 49 5 '  +1.08%  '  # E49: was '  -0.95%  '
Set-TextCell 50 2 'Cronos'  # B50: was 'MultiversX'
Set-TextCell 50 3 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'  # C50: was 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextCell 50 4 '0.0764'  # D50: was '45.35'
Set-TextCell 50 5 '  +8.36%  '  # E50: was '  -7.63%  '
Set-TextCell 51 2 'MultiversX'  # B51: was 'Cronos'
Set-TextCell 51 3 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'  # C51: was 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell 51 4 '43.83'  # D51: was '0.0760'
Set-TextCell 51 5 '  -7.40%  '  # E51: was '  +7.38%  '
